# The commit removes the "矩形 51" shape (id=52, text "更多相關" /
# "More related") that was the last shape on slide 1.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$shapes = $s.Shapes
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $shape = $shapes.Item($i)
    if ($shape.Name -eq "矩形 51") {
        $shape.Delete()
    }
}
